# semana 18 de 2025
# Adds week-18 ("18") data as new column U, corrects a value that was
# mis-attributed to week 17 (T32) by splitting it across T32/U32, backfills
# a previously-missing T40 value, and inserts a new reporting unit
# (cod_pre 6600103414 / cod_sub 01) as row 51 (pushing the last three rows
# down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert the new row for cod_pre 6600103414 / cod_sub 01 at row 51.
#    This shifts the previous rows 51-53 (EPMSC PEREIRA, SANIDAD POLICIA
#    NACIONAL RISARALDA, BATALLON SAN MATEO) down to rows 52-54.
# ------------------------------------------------------------------
$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = "'6600103414"
$ws.Range("B51").Value = "'01"
$ws.Range("U51").Value = 0

# ------------------------------------------------------------------
# 2) Fix up week-17 column T: row 32's old combined value (52) was
#    actually 26 (week 17) + 9 (week 18); row 40 had no week-17 figure
#    recorded before and is now backfilled with 0.
# ------------------------------------------------------------------
$ws.Range("T32").Value = 26
$ws.Range("T40").Value = 0

# ------------------------------------------------------------------
# 3) Add the new week-18 header in U1, matching the style (bold + centered)
#    already used by the other week-number headers in row 1.
# ------------------------------------------------------------------
$ws.Range("U1").Value = "'18"
$ws.Range("U1").Font.Bold = $true
$ws.Range("U1").HorizontalAlignment = -4108  # xlCenter

# ------------------------------------------------------------------
# 4) Populate the week-18 values down column U for every data row that
#    already carries a week-17 (T) figure.
# ------------------------------------------------------------------
$weekValues = @(
    @(2, 0),
    @(4, 0),
    @(5, 0),
    @(6, 13),
    @(7, 4),
    @(8, 24),
    @(10, 0),
    @(11, 0),
    @(12, 0),
    @(13, 0),
    @(17, 0),
    @(20, 0),
    @(22, 0),
    @(23, 0),
    @(26, 0),
    @(27, 5),
    @(28, 22),
    @(29, 0),
    @(30, 0),
    @(32, 9),
    @(33, 1),
    @(34, 0),
    @(35, 0),
    @(37, 0),
    @(38, 0),
    @(39, 0),
    @(40, 0),
    @(41, 0),
    @(42, 0),
    @(43, 0),
    @(44, 0),
    @(45, 0),
    @(46, 0),
    @(47, 0),
    @(48, 0),
    @(49, 0),
    @(50, 0),
    @(52, 0),
    @(53, 0),
    @(54, 0)
)

foreach ($pair in $weekValues) {
    $rowNum = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($rowNum, 21).Value = $val
}
